$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 38 from 45207 to 45208
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
